$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.733.64"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "1.878.18"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'332.58"
$ws.Range("E5").Value = "  +2.68%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").Value = "'0.4716"
$ws.Range("E7").Value = "  +3.84%  "
$ws.Range("D8").Value = "'0.3956"
$ws.Range("E8").Value = "  +2.05%  "
$ws.Range("D9").Value = "'47.81"
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("D10").Value = "'0.08047"
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("D11").Value = "'1.031"
$ws.Range("D12").Value = "'22.16"
$ws.Range("E12").Value = "  +3.78%  "
$ws.Range("D13").Value = "1.871.06"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("E14").Value = "  +1.50%  "
$ws.Range("D15").Value = "'7.128"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "'0.00001049"
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("D18").Value = "'87.12"
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").Value = "'0.06672"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Value = "27.756.20"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").Value = "'5.530"
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("D24").Value = "'11.02"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").Value = "'2.305"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("D26").Value = "2.096.43"
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("D27").Value = "'159.33"
$ws.Range("E27").Value = "  +3.54%  "
$ws.Range("D28").Value = "'20.18"
$ws.Range("D29").Value = "'2.104"
$ws.Range("E29").Value = "  +2.36%  "
$ws.Range("D30").Value = "'5.580"
$ws.Range("E30").Value = "  +2.24%  "
$ws.Range("D31").Value = "'121.82"
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("D32").Value = "'0.9841"
$ws.Range("E32").Value = "  +5.13%  "
$ws.Range("D33").Value = "'0.09518"
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("D34").Value = "'1.447"
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("D35").Value = "'3.599"
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("D36").Value = "'5.358"
$ws.Range("E36").Value = "  +1.97%  "
$ws.Range("D37").Value = "'0.06125"
$ws.Range("E37").Value = "  +1.89%  "
$ws.Range("D38").Value = "'0.02264"
$ws.Range("E38").Value = "  +1.78%  "
$ws.Range("D39").Value = "'1.232"
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("D40").Value = "'8.112"
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("D41").Value = "'0.6021"
$ws.Range("E41").Value = "  +1.88%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.1902"
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'10.36"
$ws.Range("E43").Value = "  +2.04%  "
$ws.Range("B44").Value = "WEMIXTOKEN"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'1.255"
$ws.Range("E44").Value = "  -2.14%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.5708"
$ws.Range("E45").Value = "  +1.88%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.21"
$ws.Range("E46").Value = "  +1.60%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'1.951"
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "'3.387"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.06914"
$ws.Range("E49").Value = "  +2.68%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'114.27"
$ws.Range("E50").Value = "  +5.48%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "'1.069"
$ws.Range("E51").Value = "  +2.01%  "
